$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> ECs)
$ws.Range("I2").Value = 0.6630574638774661
$ws.Range("J2").Value = 0.663057463877466
$ws.Range("M2").Value = 3.618510333333333
$ws.Range("N2").Value = 10.855531
$ws.Range("O2").Value = 0.1815566256530994
$ws.Range("P2").Value = 0.1815566256530994
$ws.Range("Q2").Value = 0.5898702558182222
$ws.Range("R2").Value = 5.308832302364
$ws.Range("S2").Value = 0.1203824757556946
$ws.Range("T2").Value = 0.1203824757556946

# Row 3 (FAPs -> FAPs)
$ws.Range("I3").Value = 0.6630574638774661
$ws.Range("J3").Value = 0.663057463877466
$ws.Range("O3").Value = 0.1937079481987336
$ws.Range("P3").Value = 0.1937079481987336
$ws.Range("S3").Value = 0.1284395008655599
$ws.Range("T3").Value = 0.1284395008655599

# Row 4 (FAPs -> MuSCs)
$ws.Range("I4").Value = 0.6630574638774661
$ws.Range("J4").Value = 0.663057463877466
$ws.Range("M4").Value = 12.45127566666667
$ws.Range("N4").Value = 37.353827
$ws.Range("O4").Value = 0.6247354261481669
$ws.Range("P4").Value = 0.6247354261481669
$ws.Range("Q4").Value = 2.029740552376444
$ws.Range("R4").Value = 18.267664971388
$ws.Range("S4").Value = 0.4142354872562116
$ws.Range("T4").Value = 0.4142354872562115

# Row 5 (MuSCs -> ECs)
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.08283833333333333
$ws.Range("H5").Value = 0.248515
$ws.Range("I5").Value = 0.3369425361225339
$ws.Range("J5").Value = 0.3369425361225339
$ws.Range("M5").Value = 3.618510333333333
$ws.Range("N5").Value = 10.855531
$ws.Range("O5").Value = 0.1815566256530994
$ws.Range("P5").Value = 0.1815566256530994
$ws.Range("Q5").Value = 0.2997513651627777
$ws.Range("R5").Value = 2.697762286465
$ws.Range("S5").Value = 0.06117414989740481
$ws.Range("T5").Value = 0.06117414989740481

# Row 6 (MuSCs -> FAPs)
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.08283833333333333
$ws.Range("H6").Value = 0.248515
$ws.Range("I6").Value = 0.3369425361225339
$ws.Range("J6").Value = 0.3369425361225339
$ws.Range("O6").Value = 0.1937079481987336
$ws.Range("P6").Value = 0.1937079481987336
$ws.Range("Q6").Value = 0.3198132907933334
$ws.Range("R6").Value = 2.87831961714
$ws.Range("S6").Value = 0.06526844733317375
$ws.Range("T6").Value = 0.06526844733317373

# Row 7 (MuSCs -> MuSCs)
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.08283833333333333
$ws.Range("H7").Value = 0.248515
$ws.Range("I7").Value = 0.3369425361225339
$ws.Range("J7").Value = 0.3369425361225339
$ws.Range("M7").Value = 12.45127566666667
$ws.Range("N7").Value = 37.353827
$ws.Range("O7").Value = 0.6247354261481669
$ws.Range("P7").Value = 0.6247354261481669
$ws.Range("Q7").Value = 1.031442924100555
$ws.Range("R7").Value = 9.282986316904998
$ws.Range("S7").Value = 0.2104999388919554
$ws.Range("T7").Value = 0.2104999388919553
